$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the shared text value from "testing_route" to "acma_check"
# This value is used by both B2 and B3 (they share the same string)
$ws.Range("B2").Value = "acma_check"
$ws.Range("B3").Value = "acma_check"

# Update the selection to B2:B3 with active cell B2
$ws.Activate()
$ws.Range("B2:B3").Select()
